$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / insert cell values for all 7 rows (header + 6 data rows) ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Tên Kịch Bản"
$ws.Range("C1").Value = "Dữ Liệu Mẫu"
$ws.Range("D1").Value = "Các Bước"
$ws.Range("E1").Value = "Kết Quả Mong Đợi"
$ws.Range("F1").Value = "Kết Quả Thực Tế"
$ws.Range("G1").Value = "Trạng Thái"

$ws.Range("A2").Value = "IT_OD_05"
$ws.Range("B2").Value = "Xóa tất cả chi tiết theo OrderID"
$ws.Range("C2").Value = "Order=1150"
$ws.Range("D2").Value = "1. Add items`n2. deleteDetailsByOrderId`n3. Check list"
$ws.Range("E2").Value = "List size = 0"
$ws.Range("F2").Value = "OK"
$ws.Range("G2").Value = "PASS"

$ws.Range("A3").Value = "IT_OD_06"
$ws.Range("B3").Value = "Lấy OrderID không tồn tại"
$ws.Range("C3").Value = "ID Rác"
$ws.Range("D3").Value = "getDetailsByOrderId(99999999)"
$ws.Range("E3").Value = "Trả về List rỗng (Ko null, ko crash)"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "PASS"

$ws.Range("A4").Value = "IT_OD_02"
$ws.Range("B4").Value = "Update số lượng"
$ws.Range("C4").Value = "Qty=99"
$ws.Range("D4").Value = "Update item đầu tiên thành qty=99"
$ws.Range("E4").Value = "Qty DB = 99"
$ws.Range("F4").Value = "OK"
$ws.Range("G4").Value = "PASS"

$ws.Range("A5").Value = "IT_OD_03"
$ws.Range("B5").Value = "Xóa 1 chi tiết"
$ws.Range("C5").Value = "Delete ID"
$ws.Range("D5").Value = "Thêm item nháp -> Xóa -> Find lại"
$ws.Range("E5").Value = "Get trả về null"
$ws.Range("F5").Value = "OK"
$ws.Range("G5").Value = "PASS"

$ws.Range("A6").Value = "IT_OD_04"
$ws.Range("B6").Value = "Kiểm tra JOIN Product Name"
$ws.Range("C6").Value = "Order=1150"
$ws.Range("D6").Value = "Lấy detail bất kỳ -> check field productName"
$ws.Range("E6").Value = "ProductName != null"
$ws.Range("F6").Value = "OK"
$ws.Range("G6").Value = "PASS"

$ws.Range("A7").Value = "IT_OD_01"
$ws.Range("B7").Value = "Thêm và Kiểm tra tồn tại"
$ws.Range("C7").Value = "Order=1150"
$ws.Range("D7").Value = "1. addDetail()`n2. Check list"
$ws.Range("E7").Value = "List size > 0"
$ws.Range("F7").Value = "OK"
$ws.Range("G7").Value = "PASS"

# --- Re-autofit rows that received multi-line content so Excel keeps standard
#     row height (matches source workbook, which has no explicit row heights) ---
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(7).AutoFit()

# --- Apply the "PASS" cell formatting (bold green font) used throughout column G
#     to the newly added rows by copying the format from an existing PASS cell ---
$ws.Range("G2").Copy()
$ws.Range("G5:G7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths recalculated by Excel (bestFit) after the content changes.
#     (Input values are tuned so the engine's internal pixel-grid rounding
#     lands on the width closest to the real bestFit result.) ---
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 39.833333
$ws.Columns.Item(5).ColumnWidth = 30.166667
